$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.140.71"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").Value = "2.223.84"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'293.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.25%  "

$ws.Range("D6").Value = "'87.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.51%  "

$ws.Range("E7").Value = "  -0.36%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "'0.470"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.41%  "

$ws.Range("D10").Value = "'30.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.08%  "

$ws.Range("D11").Value = "'50.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.62%  "

$ws.Range("D12").Value = "'0.0782"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").Value = "'0.114"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.35%  "

$ws.Range("E14").Value = "  -0.30%  "

$ws.Range("D15").Value = "2.582.23"
$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("E16").Value = "  -1.19%  "

$ws.Range("D17").Value = "2.245.71"
$ws.Range("E17").Value = "  +1.41%  "

$ws.Range("D18").Value = "'0.737"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.19%  "

$ws.Range("D19").Value = "40.068.24"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("E20").Value = "  +0.32%  "

$ws.Range("D21").Value = "'11.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.99%  "

$ws.Range("E22").Value = "  -0.56%  "

$ws.Range("D23").Value = "'65.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.27%  "

$ws.Range("D24").Value = "'236.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.43%  "

$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("E26").Value = "  +0.46%  "

$ws.Range("D27").Value = "'1.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "

$ws.Range("E28").Value = "  +2.57%  "

$ws.Range("E29").Value = "  +0.93%  "

$ws.Range("E30").Value = "  -6.57%  "

$ws.Range("D31").Value = "'158.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.65%  "

$ws.Range("D32").Value = "'31.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.07%  "

$ws.Range("D33").Value = "'0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "

$ws.Range("E34").Value = "  +0.36%  "

$ws.Range("D35").Value = "'3.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.63%  "

$ws.Range("E36").Value = "  -0.55%  "

$ws.Range("E37").Value = "  -2.89%  "

$ws.Range("E38").Value = "  +1.36%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.0998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.15%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'1.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.82%  "

$ws.Range("D41").Value = "'15.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.69%  "

$ws.Range("D42").Value = "2.086.74"
$ws.Range("E42").Value = "  -0.47%  "

$ws.Range("E43").Value = "  -2.49%  "

$ws.Range("D44").Value = "'19.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.56%  "

$ws.Range("D45").Value = "'10.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.37%  "

$ws.Range("E46").Value = "  +0.76%  "

$ws.Range("E47").Value = "  +2.84%  "

$ws.Range("E48").Value = "  -10.64%  "

$ws.Range("D49").Value = "2.451.35"
$ws.Range("E49").Value = "  +0.86%  "

$ws.Range("E50").Value = "  +2.21%  "

$ws.Range("E51").Value = "  +3.42%  "
